$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A ("cs") for all new rows
$ws.Cells.Item(253, 1).Value = "cs"
$ws.Cells.Item(254, 1).Value = "cs"
$ws.Cells.Item(255, 1).Value = "cs"
$ws.Cells.Item(256, 1).Value = "cs"
$ws.Cells.Item(257, 1).Value = "cs"
$ws.Cells.Item(258, 1).Value = "cs"
$ws.Cells.Item(259, 1).Value = "cs"
$ws.Cells.Item(260, 1).Value = "cs"
$ws.Cells.Item(261, 1).Value = "cs"
$ws.Cells.Item(262, 1).Value = "cs"
$ws.Cells.Item(263, 1).Value = "cs"
$ws.Cells.Item(264, 1).Value = "cs"
$ws.Cells.Item(265, 1).Value = "cs"
$ws.Cells.Item(266, 1).Value = "cs"
$ws.Cells.Item(267, 1).Value = "cs"
$ws.Cells.Item(268, 1).Value = "cs"
$ws.Cells.Item(269, 1).Value = "cs"
$ws.Cells.Item(270, 1).Value = "cs"
$ws.Cells.Item(271, 1).Value = "cs"

# Key/value cells in the exact original authoring order (preserves sharedStrings ordering)
$ws.Cells.Item(253, 2).Value = "lab.cotton.name.label"
$ws.Cells.Item(253, 3).Value = "Název"
$ws.Cells.Item(254, 2).Value = "lab.cotton.vendorId.label"
$ws.Cells.Item(254, 3).Value = "Výrobce"
$ws.Cells.Item(255, 2).Value = "lab.build.cottonId.label"
$ws.Cells.Item(255, 3).Value = "Vata"
$ws.Cells.Item(256, 2).Value = "lab.cotton.tooltip.create"
$ws.Cells.Item(256, 3).Value = "Přidat vatu"
$ws.Cells.Item(257, 2).Value = "lab.cotton.create.title"
$ws.Cells.Item(257, 3).Value = "Přidat vatu"
$ws.Cells.Item(258, 2).Value = "lab.cotton.create.subtitle"
$ws.Cells.Item(258, 3).Value = "Přidejte vatu pro použití v buildech."
$ws.Cells.Item(259, 2).Value = "lab.cotton.description.label"
$ws.Cells.Item(259, 3).Value = "Popis"
$ws.Cells.Item(260, 2).Value = "lab.cotton.create.submit"
$ws.Cells.Item(260, 3).Value = "Vytvořit vatu"
$ws.Cells.Item(261, 2).Value = "lab.cotton.create.success"
$ws.Cells.Item(261, 3).Value = "Vata byla uložena."
$ws.Cells.Item(262, 2).Value = "lab.build.coils.label"
$ws.Cells.Item(262, 3).Value = "Počet spirálek"
$ws.Cells.Item(263, 2).Value = "lab.build.coil.label"
$ws.Cells.Item(263, 3).Value = "Pozice spirálky"
$ws.Cells.Item(264, 2).Value = "lab.build.coil.label.tooltip"
$ws.Cells.Item(264, 3).Value = "Pozice spirálky je relativní umístění proti vzduchu, kdy přesné umístění je nula, umístění výše je kladné číslo a umístění níže záporné číslo. Smyslem je napovědět, jak byla spirálka umístěna a jaký byl výsledný vliv na požitek."
$ws.Cells.Item(265, 2).Value = "lab.build.cotton.label"
$ws.Cells.Item(266, 2).Value = "lab.build.cotton.label.tooltip"
$ws.Cells.Item(265, 3).Value = "Množství vaty"
$ws.Cells.Item(266, 3).Value = "Množství vaty je relativní množství, kdy nula znamená optimální usazení ve spirálce, kladná čísla znamenají větší množství vaty a záporná naopak menší množství vaty. Smyslem je napovědět, jak který atomizér snáží jaké množství vaty."
$ws.Cells.Item(267, 2).Value = "lab.build.common.title"
$ws.Cells.Item(267, 3).Value = "Obecné"
$ws.Cells.Item(268, 2).Value = "lab.build.coil.title"
$ws.Cells.Item(268, 3).Value = "Spirálka"
$ws.Cells.Item(269, 2).Value = "lab.build.cotton.title"
$ws.Cells.Item(269, 3).Value = "Vata"
$ws.Cells.Item(270, 2).Value = "lab.build.description.label"
$ws.Cells.Item(270, 3).Value = "Popis"
$ws.Cells.Item(271, 2).Value = "lab.build.ohm.label"
$ws.Cells.Item(271, 3).Value = "Odpor buildu"

# Apply the same style as the rest of the data rows (wrap text, 10pt font)
$newRange = $ws.Range("A253:C271")
$newRange.WrapText = $true
$newRange.Font.Size = 10

# Rows with long tooltip text auto-expand to 3 lines in Excel
$ws.Rows.Item(264).RowHeight = 45
$ws.Rows.Item(266).RowHeight = 45

# Update viewport to match the new scroll position / selection
$null = $ws.Range("B261").Select()
$excel.ActiveWindow.ScrollRow = 245
